$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the previously-active sheet ("Partida 8"): clear the lingering
#    number-format style on A22 (it was the only cell using that custom
#    cellXf) and move the cell selection to H23.
# ---------------------------------------------------------------------------
$prev = $wb.Worksheets.Item("Partida 8")
$prev.Range("A22").ClearFormats()
[void]$prev.Range("H23").Select()

# ---------------------------------------------------------------------------
# 2. Add the new "Partida 9" sheet after the last existing sheet.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Partida 9"

# Header row
$ws.Range("A1").Value = "T"
$ws.Range("B1").Value = "V"
$ws.Range("C1").Value = "A"

# Data rows (A2:C25)
$data = @(
    @(10, -5, 15),
    @(-5, 15, 10),
    @(-5, 20, 10),
    @(-5, 20, -5),
    @(20, -5, 15),
    @(30, -5, -5),
    @(-5, -5, 30),
    @(25, -5, 35),
    @(35, -10, 35),
    @(-5, 40, 25),
    @(30, -10, -5),
    @(35, -10, 30),
    @(35, 40, -5),
    @(30, 30, -5),
    @(25, -5, 20),
    @(20, -5, -10),
    @(25, -5, -10),
    @(15, -15, 10),
    @(15, -10, -5),
    @(-5, 15, -5),
    @(-5, 20, -5),
    @(20, -5, 10),
    @(10, -5, 10),
    @(-5, 15, 10)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# Match the saved selection/active-cell state on the new sheet.
[void]$ws.Range("K20").Select()
